$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Mark a few existing rows as DONE
$ws.Range("D3").Value = "DONE"
$ws.Range("D6").Value = "DONE"
$ws.Range("D7").Value = "DONE"
$ws.Range("D12").Value = "DONE"

# Fix typo in row 29 (Battly -> Battle)
$ws.Range("B29").Value = "End Battle Defeat"

# Add three new completed SFX entries
$ws.Range("B31").Value = "Blade Cuts Flesh"
$ws.Range("C31").Value = "SFX"
$ws.Range("D31").Value = "DONE"

$ws.Range("B32").Value = "Merrganaut hit"
$ws.Range("C32").Value = "SFX"
$ws.Range("D32").Value = "DONE"

$ws.Range("B33").Value = "unit hit (Human)"
$ws.Range("C33").Value = "SFX"
$ws.Range("D33").Value = "DONE"

# Update the active cell selection state to G8 (matches authoring session)
$ws.Range("G8").Select()
